$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty / new cells with values (row 14)
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 5

# Row 15 new cell
$ws.Range("I15").Value = 5

# Copy style from existing similarly-styled cells so the new cellXfs entries
# match the expected diff (I14 like I20/J20 green-fill bordered style,
# I15 like I9/I19 bordered style, J14 new plain-border style based on H14).
$ws.Range("I20").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I19").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H14").Copy() | Out-Null
$ws.Range("J14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# Restore view state: frozen pane top-left cell and active selection
$ws.Activate()
$ws.Range("J14").Select() | Out-Null
$av = $excel.ActiveWindow
$av.ScrollRow = 4
$av.ScrollColumn = 4
